$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the full data table (rows 2-51) to match the updated schedule,
# which now includes additional Afternoon-shift staff inserted into the
# existing groups plus new Sort B5-B9 staff appended at the end (dimension
# grows from A1:F45 to A1:F51).
$data = @(
    @("BMO", " Dale ", "Driscoll", "Day", "Location", "Station"),
    @("Bin Filler", " Karen", "Lohse", "Day", "Location", "Station"),
    @("Forklift", " George", "Dunn", "Day", "Location", "Station"),
    @("Forklift", " Don", "Coles", "Day", "Location", "Station"),
    @("Forklift", " George C", "Brown", "Day", "Location", "Station"),
    @("Forklift", " Ron", "Engene", "Afternoon", "Location", "Station"),
    @("Forklift", " Jerry", "Engene", "Afternoon", "Location", "Station"),
    @("Forklift", " Fay", "Lee", "Afternoon", "Location", "Station"),
    @("Line Operator", " Parveen", "Gopal", "Day", "Location", "Station"),
    @("Line Operator", " Lori", "Carter", "Day", "Location", "Station"),
    @("Line Operator", " Linda", "St. Amand", "Afternoon", "Location", "Station"),
    @("QC", " Isabel", "Roseen", "Day", "Location", "Station"),
    @("QC", " Wendy", "Casorso", "Day", "Location", "Station"),
    @("QC", " Shannon", "Fehr", "Afternoon", "Location", "Station"),
    @("Palletizer", " Javed", "Ali", "Day", "Location", "Station"),
    @("Stamping", " Jane", "Wu", "Day", "Location", "Station"),
    @("Non Rotational", " Elaine", "Roseen", "Day", "Location", "Station"),
    @("Non Rotational", " Janeanne", "Reiswig", "Day", "Location", "Station"),
    @("Non Rotational", " Sandra", "Martin", "Afternoon", "Location", "Station"),
    @("Non Rotational", " Joyce", "Salga", "Afternoon", "Location", "Station"),
    @("A Bliss", "Pauline", "Palatin", "Day", "Location", "Station"),
    @("B Bliss", " Cheryl", "Deboer", "Day", "Location", "Station"),
    @("F/L Operator", " Margie", "Butcher", "Day", "Location", "Station"),
    @("F/L Dumper (6:45am)", " Juanita", "Windels", "Day", "Location", "Station"),
    @("F/L Paper", " Branden", "Dubrett", "Day", "Location", "Station"),
    @("Sort - A1", " Jean", "Strachan", "Day", "Location", "Station"),
    @("Sort - A2", " Janice", "Koyama", "Day", "Location", "Station"),
    @("Sort - A3", " Paul", "Jansen", "Day", "Location", "Station"),
    @("Sort - A4", " Ed", "Fehr", "Day", "Location", "Station"),
    @("Sort - A5", " Gerald", "Kunz", "Day", "Location", "Station"),
    @("Sort - A6", " Gurdev", "Bains", "Day", "Location", "Station"),
    @("Sort - A7", " Dora", "Strachan", "Day", "Location", "Station"),
    @("Sort - A8", " Ann", "Cloutier", "Day", "Location", "Station"),
    @("Flow Control - A8", " Gail", "Johnson", "Day", "Location", "Station"),
    @("Flow Control - A9", " Nadine", "Boltz", "Day", "Location", "Station"),
    @("Flow Control - A10", " Paul", "Donnely", "Day", "Location", "Station"),
    @("Flow Control - A11", " Ross", "Izod", "Day", "Location", "Station"),
    @("Fill - A12", " Judy", "Skrove", "Day", "Location", "Station"),
    @("Fill - A13", " Doris", "Reynolds", "Day", "Location", "Station"),
    @("Fill - A14", " Tim", "Skrove", "Day", "Location", "Station"),
    @("Fill - A15", " Rick", "Sehn", "Day", "Location", "Station"),
    @("Fill - B1", " Kathy", "Yates", "Day", "Location", "Station"),
    @("Fill - B2", " Barb", "Silvester", "Day", "Location", "Station"),
    @("Fill - B3", " Judy", "Matsalla", "Day", "Location", "Station"),
    @("Fill - B4", " Debbie", "Szing", "Day", "Location", "Station"),
    @("Sort - B5", " Marcella", "Bartolomeoli", "Day", "Location", "Station"),
    @("Sort - B6", "Bob", "Young", "Day", "Location", "Station"),
    @("Sort - B7", " Jennifer ", "Knight", "Day", "Location", "Station"),
    @("Sort - B8", " Cindy", "Stubbs", "Day", "Location", "Station"),
    @("Sort - B9", " Judy Ann", "Seymour", "Day", "Location", "Station")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Host "Done writing rows"
